$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 becomes the closing/bordered row of its group (like rows 5 and 11) ---
$ws.Range("A11:E11").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

# --- Row 17: new translation entry (opening row of a new group) ---
$ws.Range("A6:E6").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 43.2

$ws.Cells.Item(17,3).Value = ' Team [team:] is up for the\nchallenge! I\''m sure of it!'
$ws.Cells.Item(17,1).Value = "SCRIPT/T01P01A/um2408.ssb"
$ws.Cells.Item(17,4).Value = ' Команде [team:] всё по\nплечу! Я уверен в этом!'
$ws.Cells.Item(17,5).Value = ' Ëïíàîäå [team:] âòæ ðï\nðìåœô! Ÿ ôâåñåî â üóïí!'
$ws.Cells.Item(17,2).Value = 88

# --- Row 18: new closing/bordered row (empty, like rows 5 and 11) ---
$ws.Range("A11:E11").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Cells.Item(18,2).ClearContents()

$excel.CutCopyMode = 0

# --- Update the active selection to match the authored state ---
$null = $ws.Range("D14").Select()
